$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.868.53"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.525.41"
$ws.Range("E3").Value = "  +2.91%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.03"
$ws.Range("E5").Value = "  +1.44%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.91"
$ws.Range("E6").Value = "  +0.42%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.527.17"
$ws.Range("E7").Value = "  +3.09%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.12%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +1.35%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.94%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -1.58%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  +2.00%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.123.78"
$ws.Range("E13").Value = "  +2.86%  "

# Row 14 - Avalanche(was ShibaInu)
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.11"
$ws.Range("E14").Value = "  +2.87%  "

# Row 15 - ShibaInu(was Avalanche)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.523.75"
$ws.Range("E16").Value = "  +3.28%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.39%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "64.872.17"
$ws.Range("E18").Value = "  +0.18%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.96"
$ws.Range("E19").Value = "  +4.56%  "

# Row 20 - Polkadot
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  +0.11%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  +4.15%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.73"
$ws.Range("E22").Value = "  +1.20%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.573"
$ws.Range("E23").Value = "  +4.47%  "

# Row 24 - WrappedeETH
$ws.Range("D24").Value = "3.664.22"
$ws.Range("E24").Value = "  +2.83%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.75"
$ws.Range("E25").Value = "  +2.28%  "

# Row 26 - Dai
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +6.45%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.60"
$ws.Range("E28").Value = "  +4.58%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.02%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +3.17%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "  +1.40%  "

# Row 32 - RenzoRestakedETH
$ws.Range("D32").Value = "3.537.14"
$ws.Range("E32").Value = "  +2.91%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.00%  "

# Row 34 - EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.71"
$ws.Range("E34").Value = "  +3.52%  "

# Row 35 - Fetch.AI
$ws.Range("E35").Value = "  +13.05%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  +1.05%  "

# Row 37 - Monero
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "170.02"
$ws.Range("E37").Value = "  +1.28%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +6.13%  "

# Row 39 - Aptos
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.82"
$ws.Range("E39").Value = "  +1.68%  "

# Row 40 - NEARProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.93"
$ws.Range("E40").Value = "  +6.72%  "

# Row 41 - Hedera
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  +5.37%  "

# Row 42 - Mantle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.819"
$ws.Range("E42").Value = "  +1.15%  "

# Row 43 - EnergySwap
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.64"
$ws.Range("E43").Value = "  +17.28%  "

# Row 44 - OKB(was FirstDigitalUSD)
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.47"
$ws.Range("E44").Value = "  +1.41%  "

# Row 45 - FirstDigitalUSD(was OKB)
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - Filecoin
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  +2.53%  "

# Row 47 - ONDO
$ws.Range("E47").Value = "  +6.37%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +2.00%  "

# Row 49 - Cosmos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.85"
$ws.Range("E49").Value = "  +6.11%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.397.44"
$ws.Range("E50").Value = "  +10.70%  "

# Row 51 - Bittensor
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "303.57"
$ws.Range("E51").Value = "  +11.17%  "
